$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Julio de 2020 a las 22:37"

# Update country data rows (name + B:H stats) per latest COVID-19 snapshot
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 3751251
$ws.Cells.Item(4, 3).Value = 56226
$ws.Cells.Item(4, 4).Value = 1700206
$ws.Cells.Item(4, 5).Value = 1909244
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 683
$ws.Cells.Item(4, 8).Value = 141801

$ws.Cells.Item(6, 1).Value = "India"
$ws.Cells.Item(6, 2).Value = 1040457
$ws.Cells.Item(6, 3).Value = 34820
$ws.Cells.Item(6, 4).Value = 654078
$ws.Cells.Item(6, 5).Value = 360094
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 676
$ws.Cells.Item(6, 8).Value = 26285

$ws.Cells.Item(19, 1).Value = "Alemania"
$ws.Cells.Item(19, 2).Value = 202337
$ws.Cells.Item(19, 3).Value = 501
$ws.Cells.Item(19, 4).Value = 186900
$ws.Cells.Item(19, 5).Value = 6277
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 3
$ws.Cells.Item(19, 8).Value = 9160

$ws.Cells.Item(27, 1).Value = "Egipto"
$ws.Cells.Item(27, 2).Value = 86474
$ws.Cells.Item(27, 3).Value = 703
$ws.Cells.Item(27, 4).Value = 27302
$ws.Cells.Item(27, 5).Value = 54984
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 68
$ws.Cells.Item(27, 8).Value = 4188

$ws.Cells.Item(78, 1).Value = "Costa Rica"
$ws.Cells.Item(78, 2).Value = 9969
$ws.Cells.Item(78, 3).Value = 423
$ws.Cells.Item(78, 4).Value = 2818
$ws.Cells.Item(78, 5).Value = 7104
$ws.Cells.Item(78, 6).Value = 0
$ws.Cells.Item(78, 7).Value = 5
$ws.Cells.Item(78, 8).Value = 47

$ws.Cells.Item(93, 1).Value = "Guinea"
$ws.Cells.Item(93, 2).Value = 6430
$ws.Cells.Item(93, 3).Value = 71
$ws.Cells.Item(93, 4).Value = 5233
$ws.Cells.Item(93, 5).Value = 1158
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 39

$ws.Cells.Item(94, 1).Value = "Gabon"
$ws.Cells.Item(94, 2).Value = 6315
$ws.Cells.Item(94, 3).Value = 194
$ws.Cells.Item(94, 4).Value = 3865
$ws.Cells.Item(94, 5).Value = 2404
$ws.Cells.Item(94, 6).Value = 0
$ws.Cells.Item(94, 7).Value = 0
$ws.Cells.Item(94, 8).Value = 46

$ws.Cells.Item(133, 1).Value = "Ruanda"
$ws.Cells.Item(133, 2).Value = 1485
$ws.Cells.Item(133, 3).Value = 12
$ws.Cells.Item(133, 4).Value = 811
$ws.Cells.Item(133, 5).Value = 670
$ws.Cells.Item(133, 6).Value = 0
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 4

$ws.Cells.Item(169, 1).Value = "Burundi"
$ws.Cells.Item(169, 2).Value = 310
$ws.Cells.Item(169, 3).Value = 7
$ws.Cells.Item(169, 4).Value = 207
$ws.Cells.Item(169, 5).Value = 102
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = 1

$ws.Cells.Item(172, 1).Value = "Eritrea"
$ws.Cells.Item(172, 2).Value = 251
$ws.Cells.Item(172, 3).Value = 0
$ws.Cells.Item(172, 4).Value = 155
$ws.Cells.Item(172, 5).Value = 96
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 0

$ws.Cells.Item(185, 1).Value = "Barbados"
$ws.Cells.Item(185, 2).Value = 104
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 91
$ws.Cells.Item(185, 5).Value = 6
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 7

$ws.Cells.Item(188, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(188, 2).Value = 79
$ws.Cells.Item(188, 3).Value = 1
$ws.Cells.Item(188, 4).Value = 63
$ws.Cells.Item(188, 5).Value = 1
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 15

$ws.Cells.Item(189, 1).Value = "Gambia"
$ws.Cells.Item(189, 2).Value = 78
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 34
$ws.Cells.Item(189, 5).Value = 41
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 3

$ws.Cells.Item(209, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(209, 2).Value = 15
$ws.Cells.Item(209, 3).Value = 4
$ws.Cells.Item(209, 4).Value = 8
$ws.Cells.Item(209, 5).Value = 7
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

$ws.Cells.Item(210, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 2).Value = 13
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 13
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0

$ws.Cells.Item(211, 1).Value = "Islas Malvinas"
$ws.Cells.Item(211, 2).Value = 13
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 13
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 0

$ws.Cells.Item(212, 1).Value = "Montserrat"
$ws.Cells.Item(212, 2).Value = 12
$ws.Cells.Item(212, 3).Value = 0
$ws.Cells.Item(212, 4).Value = 10
$ws.Cells.Item(212, 5).Value = 1
$ws.Cells.Item(212, 6).Value = 0
$ws.Cells.Item(212, 7).Value = 0
$ws.Cells.Item(212, 8).Value = 1

$ws.Cells.Item(213, 1).Value = "Santa Sede"
$ws.Cells.Item(213, 2).Value = 11
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 12
$ws.Cells.Item(213, 5).Value = 0
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 0

$ws.Cells.Item(218, 1).Value = "San Pedro y Miquelon"
$ws.Cells.Item(218, 2).Value = 4
$ws.Cells.Item(218, 3).Value = 2
$ws.Cells.Item(218, 4).Value = 1
$ws.Cells.Item(218, 5).Value = 3
$ws.Cells.Item(218, 6).Value = 0
$ws.Cells.Item(218, 7).Value = 0
$ws.Cells.Item(218, 8).Value = 0

$ws.Cells.Item(219, 1).Value = "Anguila"
$ws.Cells.Item(219, 2).Value = 3
$ws.Cells.Item(219, 3).Value = 0
$ws.Cells.Item(219, 4).Value = 3
$ws.Cells.Item(219, 5).Value = 0
$ws.Cells.Item(219, 6).Value = 0
$ws.Cells.Item(219, 7).Value = 0
$ws.Cells.Item(219, 8).Value = 0
